$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column M into the new column N, then fill in
# the "27-jun" header and the day's values for each product row.
$ws.Range("M1:M11").Copy() | Out-Null
$ws.Range("N1:N11").PasteSpecial(-4122) | Out-Null

$ws.Range("N1").Value = "27-jun"
$ws.Range("N2").Value = 14
$ws.Range("N3").Value = 16
$ws.Range("N4").Value = 7
$ws.Range("N5").Value = 13
$ws.Range("N6").Value = 12
$ws.Range("N7").Value = 19
$ws.Range("N8").Value = 14
$ws.Range("N9").Value = 16
$ws.Range("N10").Value = 11
$ws.Range("N11").Value = 6

# Remove the picture that was embedded on the sheet.
while ($ws.Shapes.Count -gt 0) {
    $ws.Shapes.Item(1).Delete() | Out-Null
}

# Match the saved selection state.
$ws.Range("N12").Select() | Out-Null
